# Update Marilith_Profits market-price / profit figures per latest Universalis pull
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 35.375
$ws.Range("I5").Value = 35.375
$ws.Range("K5").Value = 35.375
$ws.Range("M5").Value = 79.625

$ws.Range("H9").Value = 66
$ws.Range("I9").Value = 48
$ws.Range("J9").Value = 84
$ws.Range("K9").Value = 48
$ws.Range("L9").Value = 84
$ws.Range("M9").Value = 121
$ws.Range("N9").Value = -422

$ws.Range("H40").Value = 3941.3635
$ws.Range("I40").Value = 2428.3333
$ws.Range("J40").Value = 4508.75
$ws.Range("K40").Value = 2428.3333
$ws.Range("L40").Value = 4508.75
$ws.Range("M40").Value = -2253.3333
$ws.Range("N40").Value = -4858.75

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H112").Value = 2073.9
$ws.Range("I112").Value = 1600
$ws.Range("K112").Value = 4800
$ws.Range("M112").Value = -3692

$ws.Range("H135").Value = 1204.4445
$ws.Range("I135").Value = 992.25
$ws.Range("J135").Value = 1374.2
$ws.Range("K135").Value = 8930.25
$ws.Range("L135").Value = 12367.8
$ws.Range("M135").Value = -6395.25
$ws.Range("N135").Value = -17437.8

$ws.Range("H137").Value = 2137.5625
$ws.Range("I137").Value = 1535.5834
$ws.Range("K137").Value = 4606.7502
$ws.Range("M137").Value = -2056.7502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H63").Value = 3025.6428
$ws.Range("I63").Value = 2122.7144
$ws.Range("J63").Value = 3928.5715
$ws.Range("K63").Value = 2122.7144
$ws.Range("L63").Value = 3928.5715
$ws.Range("M63").Value = -1436.7144
$ws.Range("N63").Value = -5300.5715

$ws.Range("H66").Value = 3025.6428
$ws.Range("I66").Value = 2122.7144
$ws.Range("J66").Value = 3928.5715
$ws.Range("K66").Value = 10613.572
$ws.Range("L66").Value = 19642.8575
$ws.Range("M66").Value = -7181.572
$ws.Range("N66").Value = -26506.8575

$ws.Range("H88").Value = 8984.429
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 9648.5
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 9648.5
$ws.Range("M88").Value = -4594
$ws.Range("N88").Value = -10460.5

$ws.Range("H91").Value = 8984.429
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 9648.5
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 9648.5
$ws.Range("M91").Value = -3596
$ws.Range("N91").Value = -12456.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4250

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8750.625
$ws.Range("I62").Value = 9178.333000000001
$ws.Range("K62").Value = 9178.333000000001
$ws.Range("M62").Value = -8554.333000000001

$ws.Range("H65").Value = 8750.625
$ws.Range("I65").Value = 9178.333000000001
$ws.Range("K65").Value = 45891.665
$ws.Range("M65").Value = -42771.665

$ws.Range("H107").Value = 721.9
$ws.Range("J107").Value = 868.3333
$ws.Range("L107").Value = 868.3333
$ws.Range("N107").Value = -4708.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 170.5625
$ws.Range("J2").Value = 359
$ws.Range("L2").Value = 2154
$ws.Range("N2").Value = -2380

$ws.Range("H8").Value = 1292.3334
$ws.Range("I8").Value = 1292.3334
$ws.Range("K8").Value = 3877.0002
$ws.Range("M8").Value = -3738.0002

$ws.Range("H17").Value = 125
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 600
$ws.Range("M17").Value = -131
$ws.Range("N17").Value = -938

$ws.Range("H23").Value = 531.4545000000001
$ws.Range("I23").Value = 489.5
$ws.Range("K23").Value = 1468.5
$ws.Range("M23").Value = -1233.5

$ws.Range("H34").Value = 1448.75
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15168

$ws.Range("H68").Value = 856.75
$ws.Range("I68").Value = 992.3333
$ws.Range("J68").Value = 450
$ws.Range("K68").Value = 2976.9999
$ws.Range("L68").Value = 1350
$ws.Range("M68").Value = -2165.9999
$ws.Range("N68").Value = -2972

$ws.Range("H71").Value = 856.75
$ws.Range("I71").Value = 992.3333
$ws.Range("J71").Value = 450
$ws.Range("K71").Value = 8930.9997
$ws.Range("L71").Value = 4050
$ws.Range("M71").Value = -4874.9997
$ws.Range("N71").Value = -12162

$ws.Range("H97").Value = 894.6923
$ws.Range("J97").Value = 533.1
$ws.Range("L97").Value = 1599.3
$ws.Range("N97").Value = -2591.3

$ws.Range("H131").Value = 2512.4443
$ws.Range("J131").Value = 3933.3333
$ws.Range("L131").Value = 11799.9999
$ws.Range("N131").Value = -21879.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5253.4
$ws.Range("I16").Value = 5253.4
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5253.4
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -5083.4
$ws.Range("N16").ClearContents()

$ws.Range("H122").Value = 3676.25
$ws.Range("J122").Value = 3568.3333
$ws.Range("L122").Value = 10704.9999
$ws.Range("N122").Value = -15604.9999

$ws.Range("H132").Value = 30932.666
$ws.Range("I132").Value = 27110.334
$ws.Range("K132").Value = 81331.00199999999
$ws.Range("M132").Value = -78801.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 33332.5
$ws.Range("J98").Value = 33332.5
$ws.Range("L98").Value = 33332.5
$ws.Range("N98").Value = -39322.5

$ws.Range("H107").Value = 420.25
$ws.Range("I107").Value = 420.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1260.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 659.25
$ws.Range("N107").ClearContents()

$ws.Range("H136").Value = 3170.3845
$ws.Range("I136").Value = 3407.5
$ws.Range("K136").Value = 10222.5
$ws.Range("M136").Value = -7672.5
